$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (result) to hold the match_number data.
$ws.Columns("B").Insert()

# The inserted column inherits column A's formatting on the data rows; clear it so the
# new data cells are plain/unstyled like the rest of the data cells in the table.
$ws.Range("B2:B36").ClearFormats()

# Header cell gets the same (bold/centered/bordered) style as the other header cells.
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("B1").Value = "match_number"

# Fill in match numbers for rows 2-36 (Match 2 .. Match 36).
for ($row = 2; $row -le 36; $row++) {
    $ws.Cells.Item($row, 2).Value = "Match $row"
}

$excel.CutCopyMode = $false
